# Apply the "gh-pages output" update (commit 456a3b4) to both the
# "展览" sheet and the "全部类型" sheet, which share identical content.

$wb = $excel.ActiveWorkbook
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Simple "想去人数" (F column) counter bumps -----------------------
    $ws.Range("F4").Value  = 1144
    $ws.Range("F5").Value  = 47
    $ws.Range("F6").Value  = 12320
    $ws.Range("F10").Value = 439
    $ws.Range("F12").Value = 904
    $ws.Range("F13").Value = 13604
    $ws.Range("F14").Value = 13783
    $ws.Range("F16").Value = 164
    $ws.Range("F18").Value = 40
    $ws.Range("F19").Value = 1031
    $ws.Range("F21").Value = 51

    # --- Insert a new row at position 23 ------------------------------
    # This shifts the old row 23 ("Come in joy") down to row 24, and
    # row 22 keeps its place but its content is replaced below.
    $ws.Rows.Item(23).Insert()

    # Give the newly inserted row-23 cells the same look as the row
    # above it (A23 needs the bordered/centered style used by column A).
    $ws.Range("A22").Copy($ws.Range("A23"))

    # New row 23 reuses the original "OCG国潮动漫游戏嘉年华" data that used
    # to live in row 22, but the "想去人数" count changed from 4862 to 4876.
    $ws.Range("A23").Value = 22
    $ws.Range("B23").NumberFormat = "@"
    $ws.Range("B23").Value = "2024-05-04"
    $ws.Range("B23").ClearFormats()
    $ws.Range("C23").Value = "苏州·OCG国潮动漫游戏嘉年华"
    $ws.Range("D23").Value = "苏州大道东688号 苏州国际博览中心"
    $ws.Range("E23").Value = "2024.05.04 09:00-05.05 17:00"
    $ws.Range("F23").Value = 4876
    $ws.Range("G23").Value = 65
    $ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=82779"
    $ws.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg"

    # Row 22 becomes the new "大会员抢先购...阿杰内场" sub-event.
    $ws.Range("C22").Value = "【大会员抢先购】苏州·OCG国潮动漫游戏嘉年华阿杰内场"
    $ws.Range("E22").Value = "2024.05.04 09:00-05.04 17:00"
    $ws.Range("F22").Value = 2
    $ws.Range("G22").Value = 288
    $ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=82940"
    $ws.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202403/lLKmv48C1710511298160.jpeg"

    # Row 24 (originally row 23, now shifted down) needs its index and
    # "想去人数" updated.
    $ws.Range("A24").Value = 23
    $ws.Range("F24").Value = 213
}
